$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update image paths (column E) for each product row.
# Order matters: Excel's shared-string table appends new unique strings in the
# order they are first written, so we set cells in the same order as the
# original author did to reproduce the same shared-string indices.
$ws.Range("E2").Value = "/assets/img/productos/extintores/extintor-5-lbs-co2.png"
$ws.Range("E4").Value = "/assets/img/productos/extintores/Extintor-10-lbs-pqs.png"
$ws.Range("E3").Value = "/assets/img/productos/extintores/10-lbs-co2-1.png"
$ws.Range("E6").Value = "/assets/img/productos/extintores/EXTINTOR-5-LBS-PQS.png"
$ws.Range("E7").Value = "/assets/img/productos/extintores/EXTINTOR-20-LBS-PQS.png"
$ws.Range("E5").Value = "/assets/img/productos/extintores/1-11.png"

# Move the active selection to F4 as in the final saved file.
$ws.Range("F4").Select()
